$wb = $excel.ActiveWorkbook

function Set-TextCells($ws, $data) {
    foreach ($key in $data.Keys) {
        $ws.Range($key).Value = "'" + $data[$key]
        $ws.Range($key).ClearFormats()
    }
}

$wsFeatures = $wb.Worksheets.Item("Features")
$featuresData = @{
    "B2" = "0,725"
    "C2" = "0,871"
    "D2" = "0,791"
    "E2" = "0,996"
    "B3" = "0,725"
    "C3" = "0,902"
    "D3" = "0,804"
    "E3" = "0,953"
    "B4" = "0,758"
    "C4" = "0,938"
    "D4" = "0,838"
    "E4" = "0,998"
    "B5" = "0,717"
    "C5" = "0,922"
    "D5" = "0,807"
    "E5" = "0,962"
    "B6" = "0,469"
    "C6" = "0,469"
    "D6" = "0,469"
    "E6" = "1,000"
    "B7" = "0,554"
    "C7" = "0,607"
    "D7" = "0,580"
    "E7" = "1,000"
    "B8" = "0,293"
    "C8" = "0,400"
    "D8" = "0,338"
    "E8" = "0,970"
    "B9" = "0,545"
    "C9" = "0,750"
    "D9" = "0,632"
    "E9" = "1,000"
    "B10" = "0,531"
    "C10" = "0,895"
    "D10" = "0,667"
    "E10" = "0,877"
    "B11" = "0,490"
    "C11" = "0,641"
    "D11" = "0,556"
    "E11" = "0,981"
    "B12" = "0,451"
    "C12" = "0,622"
    "D12" = "0,523"
    "E12" = "0,966"
    "B13" = "0,520"
    "C13" = "0,634"
    "D13" = "0,571"
    "E13" = "0,981"
    "B14" = "0,500"
    "C14" = "0,641"
    "D14" = "0,562"
    "E14" = "0,959"
    "B15" = "0,607"
    "C15" = "0,872"
    "D15" = "0,716"
    "E15" = "1,000"
    "B16" = "0,516"
    "C16" = "0,971"
    "D16" = "0,673"
    "E16" = "0,822"
    "B17" = "0,460"
    "C17" = "0,697"
    "D17" = "0,554"
    "E17" = "1,000"
    "B18" = "0,516"
    "C18" = "0,640"
    "D18" = "0,571"
    "E18" = "0,962"
    "B19" = "0,429"
    "C19" = "0,840"
    "D19" = "0,568"
    "E19" = "0,687"
    "B20" = "0,408"
    "C20" = "0,833"
    "D20" = "0,548"
    "E20" = "0,654"
    "B21" = "0,481"
    "C21" = "0,591"
    "D21" = "0,531"
    "E21" = "0,938"
    "B22" = "0,636"
    "C22" = "0,438"
    "D22" = "0,519"
    "E22" = "1,000"
    "B23" = "0,476"
    "C23" = "0,588"
    "D23" = "0,526"
    "E23" = "0,952"
    "B24" = "0,556"
    "C24" = "0,357"
    "D24" = "0,435"
    "E24" = "0,950"
    "B25" = "0,450"
    "C25" = "0,529"
    "D25" = "0,486"
    "E25" = "0,947"
    "B26" = "0,250"
    "C26" = "0,154"
    "D26" = "0,190"
    "E26" = "1,000"
    "B27" = "0,577"
    "C27" = "0,750"
    "D27" = "0,652"
    "E27" = "1,000"
    "B28" = "0,217"
    "C28" = "0,588"
    "D28" = "0,317"
    "E28" = "1,000"
    "B29" = "0,500"
    "C29" = "0,750"
    "D29" = "0,600"
    "E29" = "0,706"
    "B30" = "0,692"
    "C30" = "0,750"
    "D30" = "0,720"
    "E30" = "1,000"
    "B31" = "0,250"
    "C31" = "0,300"
    "D31" = "0,273"
    "E31" = "0,779"
    "B32" = "0,152"
    "C32" = "0,455"
    "D32" = "0,227"
    "E32" = "1,000"
    "B33" = "0,667"
    "C33" = "0,909"
    "D33" = "0,769"
    "E33" = "1,000"
    "B34" = "0,429"
    "C34" = "0,545"
    "D34" = "0,480"
    "E34" = "1,000"
    "B35" = "0,389"
    "C35" = "0,778"
    "D35" = "0,519"
    "E35" = "1,000"
    "B36" = "1,000"
    "C36" = "0,429"
    "D36" = "0,600"
    "E36" = "0,429"
    "B37" = "0,312"
    "C37" = "0,556"
    "D37" = "0,400"
    "E37" = "1,000"
    "B38" = "0,167"
    "C38" = "0,143"
    "D38" = "0,154"
    "E38" = "1,000"
    "B39" = "0,625"
    "C39" = "0,625"
    "D39" = "0,625"
    "E39" = "1,000"
    "B40" = "0,600"
    "C40" = "1,000"
    "D40" = "0,750"
    "E40" = "0,600"
    "B41" = "0,556"
    "C41" = "0,714"
    "D41" = "0,625"
    "E41" = "1,000"
    "B43" = "0,500"
    "C43" = "0,667"
    "D43" = "0,571"
    "E43" = "1,000"
    "B44" = "0,750"
    "C44" = "0,500"
    "D44" = "0,600"
    "E44" = "0,883"
    "B45" = "0,750"
    "C45" = "0,600"
    "D45" = "0,667"
    "E45" = "0,863"
    "B46" = "0,182"
    "C46" = "0,500"
    "D46" = "0,267"
    "E46" = "1,000"
    "B47" = "0,625"
    "C47" = "0,556"
    "D47" = "0,588"
    "E47" = "1,000"
    "B48" = "0,625"
    "C48" = "0,625"
    "D48" = "0,625"
    "E48" = "1,000"
    "B49" = "1,000"
    "C49" = "0,500"
    "D49" = "0,667"
    "E49" = "0,500"
    "B50" = "1,000"
    "C50" = "0,667"
    "D50" = "0,800"
    "E50" = "0,667"
    "B51" = "0,750"
    "C51" = "0,750"
    "D51" = "0,750"
    "E51" = "0,718"
    "E53" = "0,967"
    "B55" = "0,333"
    "C55" = "0,333"
    "D55" = "0,333"
    "E55" = "1,000"
    "E56" = "1,000"
    "E60" = "0,391"
    "E61" = "0,145"
    "E65" = "1,000"
    "B69" = "1,000"
    "C69" = "0,500"
    "D69" = "0,667"
    "E69" = "0,500"
}
Set-TextCells $wsFeatures $featuresData

$wsGlobal = $wb.Worksheets.Item("Global Metrics")
$globalData = @{
    "B2" = "0,428"
    "C2" = "0,742"
    "D2" = "0,592"
    "E2" = "0,935"
}
Set-TextCells $wsGlobal $globalData

